$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.195367464011257
$ws.Range("C2").Value = 0.0298372839773424
$ws.Range("D2").Value = 0.911534437348919
$ws.Range("E2").Value = 0.0457841757765992
$ws.Range("F2").Value = 0.860338420463975
$ws.Range("G2").Value = 0.048165385864271
$ws.Range("H2").Value = 0.035429519789299
$ws.Range("I2").Value = 0.142764368438143
$ws.Range("J2").Value = 0.139300790128802
$ws.Range("K2").Value = 0.0423205974672584
$ws.Range("L2").Value = 0.946422773027384
$ws.Range("M2").Value = 0.868311866363604
$ws.Range("N2").Value = 0.0221163906627701
$ws.Range("O2").Value = 0.107731716996789
$ws.Range("P2").Value = 0.0115813399718584
$ws.Range("Q2").Value = 0.018616733412707
$ws.Range("R2").Value = 0.130389291770394
$ws.Range("S2").Value = 0.0222607064256593
$ws.Range("T2").Value = 0.980986398239348
$ws.Range("U2").Value = 0.0538297795576722
$ws.Range("V2").Value = 0.0521701482844464
$ws.Range("W2").Value = 0.223003932604539
$ws.Range("X2").Value = 0.106721506656565

# Row 3
$ws.Range("B3").Value = 0.0142511815853087
$ws.Range("C3").Value = 0.896922466356388
$ws.Range("D3").Value = 0.00321102572428474
$ws.Range("E3").Value = 0.777934119854241
$ws.Range("F3").Value = 0.00721578814446008
$ws.Range("G3").Value = 0.133419922791067
$ws.Range("H3").Value = 0.0700292239419851
$ws.Range("I3").Value = 0.00728794602590468
$ws.Range("J3").Value = 0.0208536277374896
$ws.Range("K3").Value = 0.136594869574629
$ws.Range("L3").Value = 0.0405166504311433
$ws.Range("M3").Value = 0.0104268138687448
$ws.Range("N3").Value = 0.0419237291193131
$ws.Range("O3").Value = 0.000865894577335209
$ws.Range("P3").Value = 0.00176786809539272
$ws.Range("Q3").Value = 0.0162355233250352
$ws.Range("R3").Value = 0.000757657755168308
$ws.Range("S3").Value = 0.0295847313922863
$ws.Range("T3").Value = 0.00220081538406032
$ws.Range("U3").Value = 0.0131327344229173
$ws.Range("V3").Value = 0.0953927192697622
$ws.Range("W3").Value = 0.00375220983511924
$ws.Range("X3").Value = 0.020781469856045

# Row 4
$ws.Range("B4").Value = 0.745426994263448
$ws.Range("C4").Value = 0.009488761409965
$ws.Range("D4").Value = 0.0168488653173143
$ws.Range("E4").Value = 0.0198073384565429
$ws.Range("F4").Value = 0.0284302052891727
$ws.Range("G4").Value = 0.797561063607172
$ws.Range("H4").Value = 0.0252552585056103
$ws.Range("I4").Value = 0.029548652451564
$ws.Range("J4").Value = 0.00389652559800844
$ws.Range("K4").Value = 0.00465418335317675
$ws.Range("L4").Value = 0.00660244615218097
$ws.Range("M4").Value = 0.00274199949489483
$ws.Range("N4").Value = 0.00122668398455821
$ws.Range("O4").Value = 0.000216473644333802
$ws.Range("P4").Value = 0.946350615145939
$ws.Range("Q4").Value = 0.0158025760363676
$ws.Range("R4").Value = 0.865822419453765
$ws.Range("S4").Value = 0.946278457264495
$ws.Range("T4").Value = 0.0131688133636396
$ws.Range("U4").Value = 0.924775408594004
$ws.Range("V4").Value = 0.828769347331962
$ws.Range("W4").Value = 0.75863188656781
$ws.Range("X4").Value = 0.864884366994985

# Row 5
$ws.Range("B5").Value = 0.044918281199264
$ws.Range("C5").Value = 0.0637154093155825
$ws.Range("D5").Value = 0.0682974347873146
$ws.Range("E5").Value = 0.156402208031172
$ws.Range("F5").Value = 0.104015586102392
$ws.Range("G5").Value = 0.0208536277374896
$ws.Range("H5").Value = 0.869285997763106
$ws.Range("I5").Value = 0.820362954143666
$ws.Range("J5").Value = 0.8359490565357
$ws.Range("K5").Value = 0.816430349604936
$ws.Range("L5").Value = 0.00638597250784717
$ws.Range("M5").Value = 0.118483241332034
$ws.Range("N5").Value = 0.934588880470469
$ws.Range("O5").Value = 0.891185914781542
$ws.Range("P5").Value = 0.0403001767868095
$ws.Range("Q5").Value = 0.94934516722589
$ws.Range("R5").Value = 0.00295847313922863
$ws.Range("S5").Value = 0.00184002597683732
$ws.Range("T5").Value = 0.00357181513150774
$ws.Range("U5").Value = 0.00793736695890609
$ws.Range("V5").Value = 0.0236677851138291
$ws.Range("W5").Value = 0.0143954973481979
$ws.Range("X5").Value = 0.00407692030161994
